$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new data rows for the Palo Alto, CA case study
$ws.Range("A43").Value = "Princeton"
$ws.Range("B43").Value = 814

$ws.Range("B46").Value = "Palo Alto, CA"

# Update the visible window / selection to match the saved view
$win = $excel.ActiveWindow
$win.ScrollRow = 36
$win.ScrollColumn = 1
$ws.Range("C47").Select()
